$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 93, shifting existing rows 93..185 down to 94..186.
$ws.Rows.Item(93).Insert()

# Populate the new row 93 with the new data record.
$ws.Range("A93").Value = 3
$ws.Range("B93").Value = "Femacal de La Calera"
$ws.Range("C93").Value = "Coquimbo"
$ws.Range("D93").Value2 = 44587
$ws.Range("D93").NumberFormat = $ws.Range("D94").NumberFormat
$ws.Range("E93").Value = 5
$ws.Range("F93").Value = 100112010
$ws.Range("G93").Value = "Achicoria"
$ws.Range("H93").Value = "Sin especificar"
$ws.Range("I93").Value = "Primera"
$ws.Range("J93").Value = 60
$ws.Range("K93").Value = 6000
$ws.Range("L93").Value = 6000
$ws.Range("M93").Value = 6000
$ws.Range("N93").Value = "`$/caja 16 unidades"
$ws.Range("O93").Value = "Provincia de Quillota"
$ws.Range("P93").Value = 375
$ws.Range("Q93").Value = 16
$ws.Range("R93").Value = "Hortaliza"
